# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
#
# The underlying scrape re-ordered a few fixtures that share the same
# kickoff date/time, so for each pair below the two result rows had their
# entire data payload (everything except the running id in column A, the
# Div in column C and the Date in column D) swapped with one another.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($sheet, $row1, $row2, $colStart, $colEnd)
    for ($col = $colStart; $col -le $colEnd; $col++) {
        $cell1 = $sheet.Cells.Item($row1, $col)
        $cell2 = $sheet.Cells.Item($row2, $col)
        $val1 = $cell1.Value2
        $val2 = $cell2.Value2
        $cell1.Value2 = $val2
        $cell2.Value2 = $val1
    }
}

# Columns B (id/odds id) and E..AD (HomeTeam .. PL_AhUnder) = columns 2 and 5-30.
# Column C (Div) and D (Date) are identical between each pair, so it is safe
# (and simpler) to just include the whole B..AD block (columns 2-30).

Swap-RowRange $ws 89 90 2 30
Swap-RowRange $ws 118 119 2 30
Swap-RowRange $ws 122 123 2 30
